# 19-01-25 Updated views.py (Date search bug fix)
#
# Resave-driven cleanup of 생산계획.xlsx:
#   1. The worksheet tab name "Sheet1" is lower-cased to "sheet1".
#   2. The quantity column's custom number format
#      ( _-* #,##0_-;\-* #,##0_-;_-* "-"_-;_-@_- , custom id 177 ) is
#      functionally identical to Excel's built-in "Comma [0]" format
#      (built-in id 41). Re-applying the format on the styled cells makes
#      Excel fold it onto that built-in id instead of keeping a redundant
#      custom entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet tab (Sheet1 -> sheet1)
$ws.Name = "sheet1"

# 2. Re-apply the (built-in-equivalent) number format to the styled
#    quantity cells D2:D20 so it resolves to built-in format id 41.
$ws.Range("D2:D20").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""_);_(@_)"
